# Weekly update: two new price records for potatoes (Papa) at
# "Terminal Hortofrutícola Agro Chillán" are inserted right after the
# existing row 161, pushing all subsequent rows down by two positions
# (old row 162 becomes row 164, ..., old row 267 becomes row 269).
#
# Excel's native Rows(...).Insert() does exactly that shift for us
# (values + formatting of every following row are carried down
# automatically), so all that is left to do is populate the two freshly
# inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 162 - everything from the old row
# 162 onward moves down by two rows, keeping its values/formatting.
$ws.Rows("162:163").Insert()

# New row 162: Asterix, "1a nueva(o)", registered 2022-01-13
$row162 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    44574,
    16,
    100114001,
    "Papa",
    "Asterix",
    "1a nueva(o)",
    400,
    7500,
    8000,
    7750,
    "`$/saco 25 kilos",
    "Región del Maule",
    310,
    25,
    "Hortaliza"
)

# New row 163: Asterix, "2a nueva(o)", registered 2022-01-13
$row163 = @(
    7,
    "Terminal Hortofrutícola Agro Chillán",
    "Ñuble",
    44574,
    16,
    100114001,
    "Papa",
    "Asterix",
    "2a nueva(o)",
    160,
    6500,
    7000,
    6750,
    "`$/saco 25 kilos",
    "Región del Maule",
    270,
    25,
    "Hortaliza"
)

for ($i = 0; $i -lt $row162.Length; $i++) {
    $ws.Cells.Item(162, $i + 1).Value = $row162[$i]
}

for ($i = 0; $i -lt $row163.Length; $i++) {
    $ws.Cells.Item(163, $i + 1).Value = $row163[$i]
}
